# Update the dSF (column F) values for the rows affected by the re-pull of
# data / recalculated mean, per the commit "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -4
    4  = -1
    5  = -4
    8  = -1
    15 = 3
    18 = -1
    23 = 6
    25 = -7
    26 = 1
    29 = -1
    34 = 2
    36 = 2
    45 = -1
    47 = 0
    57 = -1
    61 = -3
    63 = 1
    65 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
